$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-06-03 Tuesday" "2025-06-04 Wednesday"

Replace-Text "216÷6=36, 0" "908÷9=100, 8"
Replace-Text "618÷5=123, 3" "618÷9=68, 6"
Replace-Text "110÷8=13, 6" "443÷4=110, 3"
Replace-Text "403÷2=201, 1" "607÷5=121, 2"
Replace-Text "362÷6=60, 2" "995÷5=199, 0"
Replace-Text "546÷9=60, 6" "542÷9=60, 2"
Replace-Text "309÷8=38, 5" "324÷5=64, 4"
Replace-Text "314÷5=62, 4" "866÷5=173, 1"
Replace-Text "954÷7=136, 2" "167÷7=23, 6"
Replace-Text "465÷2=232, 1" "400÷8=50, 0"
Replace-Text "369÷9=41, 0" "649÷7=92, 5"
Replace-Text "524÷6=87, 2" "301÷8=37, 5"
Replace-Text "371÷3=123, 2" "185÷8=23, 1"
Replace-Text "441÷4=110, 1" "430÷6=71, 4"
Replace-Text "424÷9=47, 1" "553÷6=92, 1"
Replace-Text "451÷6=75, 1" "306÷6=51, 0"
Replace-Text "259÷2=129, 1" "584÷4=146, 0"
Replace-Text "586÷6=97, 4" "555÷4=138, 3"
Replace-Text "782÷7=111, 5" "733÷4=183, 1"
Replace-Text "373÷9=41, 4" "996÷4=249, 0"
Replace-Text "117÷5=23, 2" "253÷2=126, 1"
Replace-Text "179÷2=89, 1" "371÷6=61, 5"
Replace-Text "404÷2=202, 0" "938÷2=469, 0"
Replace-Text "931÷5=186, 1" "122÷2=61, 0"
Replace-Text "164÷8=20, 4" "423÷2=211, 1"
